$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unused "buildspec.yml" / ci guide row (old row 4), shifting rows up.
$ws.Rows("4").Delete()

# Rename the "Web Browser, Website, and Web Application Performance" topic to
# drop the leading "Web " (row 3, column C).
$ws.Range("C3").Value = "Browser, Website, and Web Application Performance"

# Adjust column widths for the Topic / Description / Generated Title Tag columns.
$ws.Columns("C").ColumnWidth = 51.666666666666664
$ws.Columns("D").ColumnWidth = 64.83333333333333
$ws.Columns("E").ColumnWidth = 26.166666666666668

# Update the selected cell.
$ws.Range("C3").Select() | Out-Null
